# "todays_orders.xlsx" / "Order History" sheet refresh: the feed now reports
# bills #586-#593 (was #575-#578) for the next day's batch, so every data row
# is replaced and four new rows are appended (the sheet grows from 4 to 8
# order rows, i.e. A1:I5 -> A1:I9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order data for rows 2-9 (Bill No, Date, Cashier, KOT, Price, SGST, CGST, Tax, Food Items)
$rows = @(
    @(593, 45749.22928240741, "Ajay Francis Anchan", 8, 60,  0.9,  0.9,  1.8,  "Chicken Burger (x1)"),
    @(592, 45749.22928240741, "Ajay Francis Anchan", 7, 140, 0.8,  0.8,  3.4,  "Chicken Cheese Burger (x1), Vanilla Shake (x2)"),
    @(591, 45749.22928240741, "Ajay Francis Anchan", 6, 80,  0.8,  0.8,  1.6,  "Belgian Coffee (x2)"),
    @(590, 45749.22928240741, "Ajay Francis Anchan", 5, 140, 2.2,  2.2,  2,    "Belgian Coffee (x1), Lime Juice (x1), Watermelon Juice (x1), Peri Peri Fries (x2)"),
    @(589, 45749.22928240741, "Ajay Francis Anchan", 4, 220, 3.1,  3.1,  6.2,  "Chicken Burger (x2), Mango Lassi (x1), Strawberry Lassi (x1)"),
    @(588, 45749.22928240741, "Ajay Francis Anchan", 3, 180, 1.2,  1.2,  4.2,  "Chicken Cheese Burger (x1), Vanilla Shake (x3)"),
    @(587, 45749.22928240741, "Ajay Francis Anchan", 2, 370, 5.03, 5.03, 10.05,"Vanilla Shake (x1), Oreo Shake (x1), Strawberry Lassi (x1), Butterscotch Lassi (x2), Chicken Wrap (x2)"),
    @(586, 45749.22928240741, "Ajay Francis Anchan", 1, 300, 3.6,  3.6,  9,    "Chicken Burger (x2), Chicken Cheese Burger (x1), Butterscotch Lassi (x2)")
)

# Extend the existing date-formatted style (column B, style index 1) down to
# the newly added rows first, by copying the format from B2 (copying the
# format only -- not the value -- keeps the existing style instead of
# creating a duplicate numFmt entry).
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Range("B2:B9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r++
}
